{"js": "// The document originally had Word's spelling/grammar-check markers\n// (<w:proofErr .../>) scattered through the paragraphs, which split otherwise\n// contiguous text into several <w:r> runs (one run per \"checked\" word plus\n// the runs in between). The edit simply removes every <w:proofErr> marker and\n// merges the runs it used to separate back into single runs per paragraph\n// (keeping paragraph properties, bookmarks, and the tab characters that\n// precede two of the code-sample lines).\n//\n// The most reliable way to reproduce that exactly with the Word JS API is to\n// rebuild each affected paragraph's content via `insertOoxml(..., Replace)`:\n// it swaps out the paragraph's children (runs + proofErr markers) for the\n// OOXML we provide, while the paragraph itself (and anything around it)\n// stays put.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst ns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapParagraphOoxml(innerParagraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    `<w:document ${ns}><w:body>` +\n    innerParagraphXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nfunction esc(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nfunction textRun(text) {\n  const preserve = /^\\s|\\s$|^$/.test(text) ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:r><w:t${preserve}>${esc(text)}</w:t></w:r>`;\n}\n\n// index -> replacement paragraph body (pPr + runs + anything else to keep)\nconst replacements = {\n  0:\n    \"<w:pPr><w:pStyle w:val=\\\"Ttulo1\\\"/></w:pPr>\" +\n    textRun(\"Unity\") +\n    textRun(\" \") +\n    textRun(\"Rigidbody \") +\n    textRun(\"3D\"),\n  2: textRun(\n    \"Is Kinematic: No se ve afectado por las f\u00edsicas y solo se puede mover por su componente Transform.\"\n  ),\n  3: textRun(\n    \"Velocity: mueve el objeto aplicando las f\u00edsicas, sin tanto Drag como el AddForce.\"\n  ),\n  4:\n    '<w:pPr><w:ind w:firstLine=\"708\"/></w:pPr>' + textRun(\"rb.velocity = Vector3\"),\n  5: textRun(\"position: posici\u00f3n del gameobject\"),\n  6: `<w:r><w:tab/><w:t>${esc(\"rb.position = Vector3\")}</w:t></w:r>`,\n  7: textRun(\"AddForce: a\u00f1ade una fuerza al objeto para empujarlo.\"),\n  8: `<w:r><w:tab/><w:t>${esc(\"rb.AddFroce(Vector3)\")}</w:t></w:r>`,\n  9:\n    \"<w:pPr><w:rPr><w:u w:val=\\\"single\\\"/></w:rPr></w:pPr>\" +\n    textRun(\n      \"infinity: cuando le ponemos en el drag o en el angular Drag infinuty, significa que el objeto para de moverse inmediatamente.\"\n    ) +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>',\n};\n\nfor (const [idx, inner] of Object.entries(replacements)) {\n  const p = paragraphs.items[Number(idx)];\n  p.insertOoxml(wrapParagraphOoxml(`<w:p>${inner}</w:p>`), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document had Word's spell/grammar-check markers (<w:proofErr/>)\n# scattered across the paragraphs; those markers split otherwise contiguous\n# sentences into several runs (one run per \"checked\" word, plus the runs of\n# text in between). The edit removes every <w:proofErr> and merges the runs\n# it used to separate back into a single run per paragraph (keeping each\n# paragraph's own properties/bookmarks, and the leading tab characters on the\n# two code-sample lines).\n#\n# `Range.InsertXML` replaces a range's contents with the OOXML fragment we\n# give it (paragraph mark, formatting and anything outside the range are left\n# alone), which is exactly what's needed to rebuild the merged runs without\n# the leftover proofErr markers.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphXml($paragraphIndex, $innerXml) {\n    $p = $d.Paragraphs.Item($paragraphIndex)\n    $xml = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" + $innerXml + \"</w:p>\"\n    $p.Range.InsertXML($xml)\n}\n\nfunction Esc($s) {\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\nfunction TextRun($text) {\n    $needsPreserve = ($text -match '^\\s') -or ($text -match '\\s$') -or ($text -eq \"\")\n    $escaped = Esc $text\n    if ($needsPreserve) {\n        return \"<w:r><w:t xml:space='preserve'>\" + $escaped + \"</w:t></w:r>\"\n    } else {\n        return \"<w:r><w:t>\" + $escaped + \"</w:t></w:r>\"\n    }\n}\n\n# Paragraph 1 (title): \"Unity\" / \" \" / \"Rigidbody \" / \"3D\"\n$r1 = TextRun \"Unity\"\n$r2 = TextRun \" \"\n$r3 = TextRun \"Rigidbody \"\n$r4 = TextRun \"3D\"\n$inner1 = \"<w:pPr><w:pStyle w:val='Ttulo1'/></w:pPr>\" + $r1 + $r2 + $r3 + $r4\nSet-ParagraphXml 1 $inner1\n\n# Paragraph 2 is empty and untouched by the diff.\n\n$inner3 = TextRun \"Is Kinematic: No se ve afectado por las f\u00edsicas y solo se puede mover por su componente Transform.\"\nSet-ParagraphXml 3 $inner3\n\n$inner4 = TextRun \"Velocity: mueve el objeto aplicando las f\u00edsicas, sin tanto Drag como el AddForce.\"\nSet-ParagraphXml 4 $inner4\n\n$r5 = TextRun \"rb.velocity = Vector3\"\n$inner5 = \"<w:pPr><w:ind w:firstLine='708'/></w:pPr>\" + $r5\nSet-ParagraphXml 5 $inner5\n\n$inner6 = TextRun \"position: posici\u00f3n del gameobject\"\nSet-ParagraphXml 6 $inner6\n\n$inner7 = \"<w:r><w:tab/><w:t>rb.position = Vector3</w:t></w:r>\"\nSet-ParagraphXml 7 $inner7\n\n$inner8 = TextRun \"AddForce: a\u00f1ade una fuerza al objeto para empujarlo.\"\nSet-ParagraphXml 8 $inner8\n\n$inner9 = \"<w:r><w:tab/><w:t>rb.AddFroce(Vector3)</w:t></w:r>\"\nSet-ParagraphXml 9 $inner9\n\n$r10 = TextRun \"infinity: cuando le ponemos en el drag o en el angular Drag infinuty, significa que el objeto para de moverse inmediatamente.\"\n$inner10 = \"<w:pPr><w:rPr><w:u w:val='single'/></w:rPr></w:pPr>\" + $r10 + \"<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>\"\nSet-ParagraphXml 10 $inner10\n"}
